# Loop working with 3 different data sets
# Update the newTours test-data sheet so that rows 2, 3 and 4 each hold a
# distinct record (instead of rows 3 and 4 duplicating row 2), and move the
# active selection to M4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("newTours")

# ----- Row 2 (data set 1) -----
$ws.Range("D2").Value  = "Acapulco"
$ws.Range("G2").Value  = "Zurich"
$ws.Range("J2").Value  = "Coach"
$ws.Range("K2").Value  = "Blue Skies Airlines"
$ws.Range("N2").Value  = "Firstname11"
$ws.Range("O2").Value  = "Lastname11"
$ws.Range("Q2").Value  = "Firstname21"
$ws.Range("R2").Value  = "Lastname21"
$ws.Range("T2").Value  = "Carte Blanche"
$ws.Range("V2").Value  = "01"
$ws.Range("AB2").Value = "Billing Address Complement 1"
$ws.Range("AG2").Value = "1, Delivery Address"
$ws.Range("AH2").Value = "1, Delivery Address Complement"

# ----- Row 3 (data set 2) -----
$ws.Range("D3").Value  = "Frankfurt"
$ws.Range("E3").Value  = "October"
$ws.Range("G3").Value  = "Sydney"
$ws.Range("H3").Value  = "November"
$ws.Range("J3").Value  = "Business"
$ws.Range("L3").Value  = "2"
$ws.Range("M3").Value  = "3"
$ws.Range("N3").Value  = "Firstname12"
$ws.Range("O3").Value  = "Lastname12"
$ws.Range("P3").Value  = "Bland"
$ws.Range("Q3").Value  = "Firstname22"
$ws.Range("R3").Value  = "Lastname22"
$ws.Range("S3").Value  = "Diabetic"
$ws.Range("U3").Value  = "9238483848"
$ws.Range("V3").Value  = "05"
$ws.Range("W3").Value  = "2005"
$ws.Range("X3").Value  = "Bob"
$ws.Range("Y3").Value  = "Gardner"
$ws.Range("Z3").Value  = "III"
$ws.Range("AA3").Value = "2, Billing Address"
$ws.Range("AB3").Value = "Billing Address Complement 2"
$ws.Range("AE3").Value = "34142"
$ws.Range("AK3").Value = "36563"
$ws.Range("AL3").Value = "UNITED STATES"

# ----- Row 4 (data set 3) -----
$ws.Range("E4").Value  = "November"
$ws.Range("H4").Value  = "December"
$ws.Range("K4").Value  = "Pangea Airlines"
$ws.Range("L4").Value  = "3"
$ws.Range("M4").Value  = "0"
$ws.Range("N4").Value  = "Firstname13"
$ws.Range("O4").Value  = "Lastname13"
$ws.Range("P4").Value  = "Kosher"
$ws.Range("Q4").Value  = "Firstname23"
$ws.Range("R4").Value  = "Lastname23"
$ws.Range("S4").Value  = "Low Sodium"
$ws.Range("T4").Value  = "American Express"
$ws.Range("U4").Value  = "4373775834"
$ws.Range("W4").Value  = "2001"
$ws.Range("X4").Value  = "Mark"
$ws.Range("Y4").Value  = "Richard"
$ws.Range("Z4").Value  = "Smith"
$ws.Range("AA4").Value = "3, Billing Address"
$ws.Range("AB4").Value = "Billing Address Complement 3"
$ws.Range("AE4").Value = "43422"
$ws.Range("AG4").Value = "3, Delivery Address"
$ws.Range("AH4").Value = "3, Delivery Address Complement"
$ws.Range("AK4").Value = "65878"
$ws.Range("AL4").Value = "BRAZIL"

# Move the active selection to M4, matching the saved sheet view state.
$ws.Range("M4").Select()

# delCountry now holds a wider value ("UNITED STATES"); let Excel re-fit the column.
$ws.Columns.Item(38).AutoFit()
